$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 64
$ws1.Range("F4").Value = 247
$ws1.Range("F6").Value = 234
$ws1.Range("F7").Value = 193
$ws1.Range("F8").Value = 1856
$ws1.Range("F9").Value = 334
$ws1.Range("F10").Value = 4320
$ws1.Range("F11").Value = 58
$ws1.Range("F12").Value = 303

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2
$ws4.Range("F4").Value = 64
$ws4.Range("F6").Value = 247
$ws4.Range("F8").Value = 234
$ws4.Range("F9").Value = 193
$ws4.Range("F12").Value = 1856
$ws4.Range("F13").Value = 334
$ws4.Range("F14").Value = 4320
$ws4.Range("F15").Value = 58
$ws4.Range("F16").Value = 303
